# 22072023 excel file write
#
# Adds a new trailing column to the "neworder", "order" and "register"
# sheets (age / country / message+registration success respectively),
# copying the formatting of an existing column and filling in the new
# header / data text.  "oldorder" only gets its active-cell selection
# nudged over to column F.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "neworder": add column F "age"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("neworder")
$ws1.Activate()

$ws1.Range("E1").Copy() | Out-Null
$ws1.Range("F1").PasteSpecial(-4122) | Out-Null
$ws1.Range("F1").Value = "age"

$ws1.Range("E2:E6").Copy() | Out-Null
$ws1.Range("F2").PasteSpecial(-4122) | Out-Null

$ws1.Range("F6").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "oldorder": only the active-cell selection moves to F1
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("oldorder")
$ws2.Activate()
$ws2.Range("F1").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "order": add column F "country"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("order")
$ws3.Activate()

$ws3.Range("E1").Copy() | Out-Null
$ws3.Range("F1").PasteSpecial(-4122) | Out-Null
$ws3.Range("F1").Value = "country"

$ws3.Range("E2:E6").Copy() | Out-Null
$ws3.Range("F2").PasteSpecial(-4122) | Out-Null

$ws3.Range("F6").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "register": add column K "message" / "registration success"
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("register")
$ws4.Activate()

$ws4.Columns.Item(11).ColumnWidth = 23.333333333333332

$ws4.Range("A1").Copy() | Out-Null
$ws4.Range("K1").PasteSpecial(-4122) | Out-Null
$ws4.Range("K1").Value = "message"

$ws4.Range("A2").Copy() | Out-Null
$ws4.Range("K2").PasteSpecial(-4122) | Out-Null
$ws4.Range("K2").Value = "registration success"

$ws4.Range("A3").Copy() | Out-Null
$ws4.Range("K3").PasteSpecial(-4122) | Out-Null

$ws4.Range("K1").Select() | Out-Null
$ws4.Activate()

Write-Output "edit complete"
